$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its original text data type —
# several new values are valid numeric literals that Excel would
# otherwise auto-convert to numbers on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.560.76'
$ws.Range("E2").Value = '  -2.47%  '
$ws.Range("D3").Value = '1.656.67'
$ws.Range("E3").Value = '  -4.19%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '213.99'
$ws.Range("E5").Value = '  -2.39%  '
$ws.Range("D6").Value = '0.510'
$ws.Range("E6").Value = '  -2.33%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '23.96'
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").Value = '0.0878'
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("D12").Value = '1.890.26'
$ws.Range("E12").Value = '  -4.25%  '
$ws.Range("D13").Value = '1.655.62'
$ws.Range("E13").Value = '  -4.32%  '
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("E16").Value = '  -2.85%  '
$ws.Range("D17").Value = '27.542.12'
$ws.Range("E17").Value = '  -2.56%  '
$ws.Range("D18").Value = '240.22'
$ws.Range("E18").Value = '  -2.12%  '
$ws.Range("E19").Value = '  -3.27%  '
$ws.Range("E20").Value = '  -4.62%  '
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  -3.95%  '
$ws.Range("D23").Value = '9.40'
$ws.Range("E23").Value = '  -2.90%  '
$ws.Range("E24").Value = '  -2.32%  '
$ws.Range("D25").Value = '145.71'
$ws.Range("D26").Value = '7.19'
$ws.Range("E26").Value = '  -4.12%  '
$ws.Range("E27").Value = '  -2.46%  '
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E31").Value = '  -2.95%  '
$ws.Range("E32").Value = '  -3.08%  '
$ws.Range("D33").Value = '1.445.93'
$ws.Range("E33").Value = '  -2.75%  '
$ws.Range("E34").Value = '  -5.09%  '
$ws.Range("E35").Value = '  -5.13%  '
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("D37").Value = '0.923'
$ws.Range("E37").Value = '  -5.46%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0171'
$ws.Range("E38").Value = '  -2.87%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.569'
$ws.Range("E39").Value = '  -5.56%  '
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("D41").Value = '69.03'
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.41'
$ws.Range("E43").Value = '  -4.37%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '0.795'
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = '2.22'
$ws.Range("E45").Value = '  -3.22%  '
$ws.Range("D46").Value = '1.798.92'
$ws.Range("E46").Value = '  -4.17%  '
$ws.Range("D47").Value = '1.70'
$ws.Range("E47").Value = '  -1.29%  '
$ws.Range("D48").Value = '88.32'
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("E49").Value = '  -6.98%  '
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("D51").Value = '7.84'
$ws.Range("E51").Value = '  -4.41%  '
